$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://www.linkedin.com/in/rogério-artigiani-cspo-psm-b882a911"
$ws.Range("B3").Value = "https://www.linkedin.com/in/cristiane-antunes"
$ws.Range("B4").Value = "https://www.linkedin.com/in/welder-soares-leitão-83a14710a"
$ws.Range("B5").Value = "https://www.linkedin.com/in/geovanaferreira"
$ws.Range("B6").Value = "https://www.linkedin.com/in/jcampos3"
$ws.Range("B7").Value = "https://www.linkedin.com/in/aline-s-bb5938117"
$ws.Range("B8").Value = "https://www.linkedin.com/in/fernanda-liveri-19a59522"
$ws.Range("B9").Value = "https://www.linkedin.com/in/anapbatista"
$ws.Range("B10").Value = "https://www.linkedin.com/in/gabrielssian"
$ws.Range("B11").Value = "https://www.linkedin.com/in/kelly-cristina-deniz-2598a722"
$ws.Range("B12").Value = "https://www.linkedin.com/in/gustavosilveira12"
$ws.Range("B13").Value = "https://www.linkedin.com/in/jéssica-cristine-alexandrino-marcondes-68232751"
$ws.Range("B14").Value = "https://www.linkedin.com/in/renan-t-ramos-1602b669"
$ws.Range("B15").Value = "https://www.linkedin.com/in/brunoac19"
$ws.Range("B16").Value = "https://www.linkedin.com/in/felipe-artur-rodrigues-batista-de-oliveira-366b3234"
$ws.Range("B17").Value = "https://www.linkedin.com/in/fabio-senatore"
$ws.Range("B18").Value = "https://www.linkedin.com/in/juliana-fava-canabrava-b44893168"
$ws.Range("B19").Value = "https://www.linkedin.com/in/isaac-de-freitas-lima-2a3567196"
$ws.Range("B20").Value = "https://www.linkedin.com/in/giangagliardo"
$ws.Range("B21").Value = "https://www.linkedin.com/in/bruno-medeiros-50396a48"
$ws.Range("B22").Value = "https://www.linkedin.com/in/markosmadeira"
$ws.Range("B23").Value = "https://www.linkedin.com/in/olívia-longarço-6a747b"
$ws.Range("B24").Value = "https://www.linkedin.com/in/guilherme-vieira-15217565"
$ws.Range("B25").Value = "https://www.linkedin.com/in/adriel-dantas-44162873"
$ws.Range("B26").Value = "https://www.linkedin.com/in/viniciusdiaspeixoto"
$ws.Range("B27").Value = "https://www.linkedin.com/in/alexsandro-souza-xavier-b5605959"
$ws.Range("B28").Value = "https://www.linkedin.com/in/thaisehagge"
$ws.Range("B29").Value = "https://www.linkedin.com/in/roberta-rodrigues-muoio-b53a4733"
$ws.Range("B30").Value = "https://www.linkedin.com/in/wesleinunes"
$ws.Range("B31").Value = "https://www.linkedin.com/in/sérgio-sobrosa-batista-30a9b61a6"
$ws.Range("B32").Value = "https://www.linkedin.com/in/juliana-carsoni-19376553"
$ws.Range("B33").Value = "https://www.linkedin.com/in/miltonbarros"
$ws.Range("B34").Value = "https://www.linkedin.com/in/tiago-vian-47a6ab11a"
$ws.Range("B35").Value = "https://www.linkedin.com/in/gustavo-p-barros-06249a67"
$ws.Range("B36").Value = "https://www.linkedin.com/in/camilasobral"
$ws.Range("B37").Value = "https://www.linkedin.com/in/álvaro-souza-a07a21a4"
$ws.Range("B38").Value = "https://www.linkedin.com/in/thiago-mendes-do-nascimento"
$ws.Range("B39").Value = "https://www.linkedin.com/in/carolinapirmez"
$ws.Range("B40").Value = "https://www.linkedin.com/in/meire-hellen-galinari-lopes-bb2a18a0"
$ws.Range("B41").Value = "https://www.linkedin.com/in/aline-miura-b8187296"
$ws.Range("B42").Value = "https://www.linkedin.com/in/daysy-andrade-silva-756876b2"
$ws.Range("B43").Value = "https://www.linkedin.com/in/andressa-alves-3703a453"
$ws.Range("B44").Value = "https://www.linkedin.com/in/thais-vergueiro-6536b926"
$ws.Range("B45").Value = "https://www.linkedin.com/in/cristina-baik-8b8b7749"
$ws.Range("B46").Value = "https://www.linkedin.com/in/patriciagoia"
$ws.Range("B47").Value = "https://www.linkedin.com/in/rayssa-albuquerque-62581539"
$ws.Range("B48").Value = "https://www.linkedin.com/in/aline-portella-17755764"
$ws.Range("B49").Value = "https://www.linkedin.com/in/silvana-fernandes-585b74"
$ws.Range("B50").Value = "https://www.linkedin.com/in/beatrizcapistrano"
$ws.Range("B51").Value = "https://www.linkedin.com/in/ronaldoschulze"
